$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text that often happens to look like a
# number ("378.13", "1.00", ...). Assigning such a string straight to .Value
# lets Excel auto-detect it as a number, which changes the stored cell type.
# Prefixing with a leading apostrophe forces text entry (as typing '378.13 at
# the keyboard would); the apostrophe itself is not stored. Excel tags the
# cell with a "Text" number format as a side effect of the apostrophe entry,
# so re-apply the Normal style right after to keep formatting untouched.
function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-PriceText "D2" '51.528.54'
$ws.Range("E2").Value = '  +1.16%  '
Set-PriceText "D3" '3.014.36'
$ws.Range("E3").Value = '  +2.31%  '
Set-PriceText "D5" '378.13'
$ws.Range("E5").Value = '  -0.26%  '
Set-PriceText "D6" '102.92'
$ws.Range("E6").Value = '  +2.04%  '
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-PriceText "D9" '0.595'
$ws.Range("E9").Value = '  +2.55%  '
Set-PriceText "D10" '36.56'
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("E11").Value = '  -0.43%  '
Set-PriceText "D13" '3.497.49'
$ws.Range("E13").Value = '  +2.83%  '
Set-PriceText "D14" '18.45'
$ws.Range("E14").Value = '  +1.28%  '
Set-PriceText "D15" '7.73'
$ws.Range("E15").Value = '  +1.43%  '
Set-PriceText "D16" '3.026.62'
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("E17").Value = '  -1.73%  '
Set-PriceText "D18" '10.52'
$ws.Range("E18").Value = '  -12.94%  '
Set-PriceText "D19" '51.533.74'
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("E20").Value = '  -0.22%  '
Set-PriceText "D21" '12.47'
$ws.Range("E21").Value = '  +0.57%  '
Set-PriceText "D22" '0.0₃0960'
$ws.Range("E22").Value = '  +1.13%  '
Set-PriceText "D23" '69.85'
$ws.Range("E23").Value = '  +0.49%  '
Set-PriceText "D24" '267.18'
$ws.Range("E24").Value = '  +0.10%  '
Set-PriceText "D25" '3.13'
$ws.Range("E25").Value = '  -3.96%  '
Set-PriceText "D26" '8.20'
$ws.Range("E26").Value = '  +1.48%  '
Set-PriceText "D27" '7.52'
$ws.Range("E27").Value = '  +5.83%  '
Set-PriceText "D28" '0.171'
$ws.Range("E28").Value = '  +5.44%  '
Set-PriceText "D29" '1.00'
$ws.Range("E29").Value = '  -0.02%  '
Set-PriceText "D30" '26.20'
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("E32").Value = '  +2.54%  '
Set-PriceText "D33" '34.11'
$ws.Range("E33").Value = '  +1.79%  '
Set-PriceText "D34" '50.68'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("E35").Value = '  +5.32%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +6.08%  '
Set-PriceText "D39" '17.19'
$ws.Range("E39").Value = '  +3.75%  '
$ws.Range("E40").Value = '  +10.22%  '
Set-PriceText "D41" '2.59'
$ws.Range("E41").Value = '  +3.71%  '
Set-PriceText "D42" '1.85'
$ws.Range("E42").Value = '  +2.17%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("E44").Value = '  +8.87%  '
Set-PriceText "D45" '122.98'
$ws.Range("E45").Value = '  +2.57%  '
Set-PriceText "D46" '21.46'
$ws.Range("E46").Value = '  +0.97%  '
Set-PriceText "D47" '2.06'
$ws.Range("E47").Value = '  +2.69%  '
Set-PriceText "D48" '2.38'
$ws.Range("E48").Value = '  +1.93%  '
Set-PriceText "D49" '2.030.80'
$ws.Range("E49").Value = '  +0.96%  '
Set-PriceText "D50" '3.314.04'
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("E51").Value = '  +1.35%  '
